# The workbook's small roster table (A1:E5) had its "Name" column (B2:B5)
# re-keyed to a different set of names. Everything else (headers, amounts,
# status, type) stays the same - only the shared-string values referenced by
# column B change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Zakeeramma"
$ws.Range("B3").Value = "Soma Mahadevamma"
$ws.Range("B4").Value = "lakshmi"
$ws.Range("B5").Value = "radha"

# Column B grew a bit wider to fit the new (longer) names and lost its
# "best fit" auto-sizing flag, becoming an explicit custom width.
$ws.Columns("B:B").ColumnWidth = 17.65

# The user's selection ended up on B6 (just below the last data row) after
# making the edits.
$ws.Range("B6").Select() | Out-Null
